# Change "mean" to "averaged" in the rolling-aggregation parameter names
# located in column N, rows 16-33 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 16; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 14)  # column N = 14
    $text = $cell.Value2
    if ($text -ne $null -and $text -like "*mean*") {
        $cell.Value = $text -replace "mean", "averaged"
    }
}

# Reflect the selected cell recorded in the saved workbook
$ws.Range("N27").Select()
